$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): F3 7560->7561, F4 285->286, F7 4188->4189
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 7561
$ws1.Range("F4").Value = 286
$ws1.Range("F7").Value = 4189

# Sheet "全部类型" (sheet4): F4 7560->7561, F6 285->286, F9 4188->4189
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 7561
$ws4.Range("F6").Value = 286
$ws4.Range("F9").Value = 4189
